$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B32 becomes a true numeric value (was stored as text "3")
$ws.Range("B32").Value = 3

# Add new row 33 with the follow-up annotation data
$ws.Range("A33").Value = "Ruilin"
$ws.Range("B33").Value = "'3"
$ws.Range("B33").ClearFormats()
$ws.Range("C33").Value = "无"
$ws.Range("D33").Value = "FBK"
$ws.Range("E33").Value = "OTH"
$ws.Range("F33").Value = "d4ad31e6-de82-4ee8-af90-c18d97ed2c36"
$ws.Range("G33").Value = "Bk7wvW-C-_annotated.xlsx"
$ws.Range("H33").Value = "We will update our paper very soon."
